$d = $word.ActiveDocument

# 1. Update the "Recommend Archaeology Condition..." heading text
$d.Content.Find.Execute("Recommend Archaeology Condition to be Satisfied", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Recommend Archaeology Condition(s)", 2) | Out-Null

# 2. Move the _GoBack bookmark from before "<Proposal Description>" to just
#    before "National Planning Policy Framework" inside the GLAAS paragraph,
#    which splits that run's text in two around the new bookmark location.
$r = $d.Content
$r.Find.Execute("National Planning Policy Framework (NPPF) and the GLAAS Charter.", $true, $false, $false, $false, $false,
                $true, 1, $false, "", 0) | Out-Null
$target = $d.Range($r.Start, $r.Start)
$d.Bookmarks.Add("_GoBack", $target) | Out-Null
